$d = $word.ActiveDocument

# Change 1: Ativação date
$d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2)

# Change 2: Programa resumido paragraph
$old2 = "Fundamentos da corrosão (termodinâmica e cinética); Principais tipos de corrosão; Controle e proteção contra a corrosão; Degradação de materiais poliméricos e cerâmicos; Oxidação em altas temperaturas."
$new2 = "Fundamentos da corrosão (termodinâmica e cinética); principais tipos de corrosão; Controle e proteção contra a corrosão; Degradação de materiais poliméricos e cerâmicos; Oxidação em altas temperaturas. Estudo de Caso"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# Change 3: Programa (full syllabus) paragraph
$old3 = "1. Corrosão e sua importância econômica, social e ambiental. 2. Aspectos termodinâmicos e cinéticos da corrosão. Polarização. 3. Diagrama de Pourbaix. 4. Tipos de corrosão (corrosão uniforme, corrosão por pites, corrosão intergranular, corrosão associado com fatores mecânicos, corrosão galvânica e corrosão atmosférica). 5. Proteção catódica e anódica. 6. Revestimentos metálicos e orgânicos. 7. Inibidores de corrosão. 8. Degradação de materiais poliméricos e cerâmicos. 9. Oxidação em altas temperaturas."
$new3 = "Corrosão e sua importância econômica, social e ambiental. 2. Aspectos termodinâmicos e cinéticos da corrosão. Polarização. Princípios básicos de eletroquímica para compreensão do fenômeno de corrosão. Diagrama de Pourbaix. 3. Tipos de corrosão (corrosão uniforme, corrosão por pites, corrosão intergranular, corrosão associado com fatores mecânicos, corrosão galvânica e corrosão atmosférica). 4. Proteção catódica e anódica. 5.Revestimentos metálicos e orgânicos. 6. Inibidores de corrosão. 7. Degradação de materiais poliméricos. 8.Degradação de materiais cerâmicos. 9.Oxidação em altas temperaturas. 10. Estudos de Casos referentes às falhas causadas em equipamentos devido processo corrosivo, atividade que potencializa o papel do aluno como protagonista do processo de ensino e aprendizagem, colocando-o em contato com problemas reais."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# Change 4: Método paragraph
$old4 = "As avaliações serão por meio de provas individuais ou trabalhos em equipes, conforme adequação ao conteúdo programático."
$new4 = "Os alunos serão avaliados continuamente quanto às habilidades gerais em função da participação ativa nas aulas. Também, serão aplicadas: provas escritas, trabalhos extraclasse, pequenos seminários e atividades para discussão dos Estudos de Caso."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)

# Change 5: Critério paragraph
$old5 = "A nota final é a média aritmética das avaliações realizadas"
$new5 = "Será considerada a média das avaliações gerais, com peso 1 (AG) somada à nota de duas provass escritas, P1 e P2, com peso 1, cada). A nota final (NF) será calculada pela equação: NF = (AG + P1+ P2)/3."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)

# Change 6: Bibliografia paragraph
$old6 = "BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Coimbra, Livraria Medina, 1996.FONTANA, M.G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro, Ed. LTC, 2007 JAMBO, H.C.M. e Fofano S. Corrosão: Fundamentos, Monitoração e Controle. Editora Ciência Moderna, 2009.JONES, D.A. Principles and Prevention of Corrosion. 2ª  Edição, Prentice Hall, 1996.McCAULEY, R. A. Corrosion of Ceramic and Composite Materials. 2ª  Edição, 2004.RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus, 1990.SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2, 2000.WEST, J. M. Basic Corrosion and Oxidation. Second Edition. New York. Ellis Horwood Limited; John Wiley & Sons. 1986."
$new6 = "ASM Handbook Committee, ASM Handbook: Corrosion: Fundamentals, Testing, and Protection, vol. 13A, Materials Park, Ohio: ASM International, 2003.BRETT, A.M.O., BRETT, C.M. Electroquímica: Princípios, métodos e aplicações. Coimbra, Livraria Medina, 1996. FONTANA, M.G. Corrosion Engineering. 3ª Edição. McGraw-Hill, 1987  GENTIL, V. Corrosão. 5ª Edição, Rio de Janeiro,Ed. LTC, 2007 JAMBO, H.C.M. e Fofano S. Corrosão: Fundamentos, Monitoração e Controle. Editora Ciência Moderna,2009. JONES, D.A. Principles and Prevention of Corrosion. 2ª Edição, Prentice Hall, 1996. McCAULEY, R. A. Corrosion of Ceramic and Composite Materials. 2ª Edição, 2004. RAMANHATAN, L. Corrosão e seu Controle. São Paulo. Ed. Hemus,1990. SHREIR, L.L., JARMAN, R.A., BURSTEIN, G.T. Corrosion. 3ª Edição. Oxford, Butterworth Heinemann, volume 2,2000. WEST, J. M. Basic Corrosion and Oxidation. Second Edition. New York. Ellis Horwood Limited; John Wiley & Sons.1986.ROBERGE,Pierre R. Corrosion Engineering Principles and PracticeMcGraw-Hill, 2008.."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)

Write-Output "Done"
